$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the "articletitle" column (C1:C21).
# This removes the shared strings for the header and the 20 article titles,
# leaving styled-but-empty cells where a style was previously applied
# (C1, C12, C13) and removing the cell entirely elsewhere.
$ws.Range("C1:C21").ClearContents()

# Move the view/selection so the sheet scrolls back to the top and the
# active cell becomes D2 (matching the saved view state in the workbook).
$ws.Range("D2").Select()
